# notes/ToDo.docx - add the 06/11/2018 status update beneath the existing
# "01/11/18" date stamp (which currently just holds the lone date line plus
# the floating "_GoBack" bookmark at the very end of the document).
#
# Approach: the new content is a run of seven fresh paragraphs (a bold+
# underlined date heading, a "Bríd:" line carrying the relocated bookmark,
# four bullet items continuing the existing "ListParagraph"/numId=1 list,
# and a closing remark paragraph). We build that slice of WordprocessingML
# by hand and drop it in with Range.InsertXML, which is the COM-exposed
# equivalent of pasting/merging OOXML into the document at a given Range.

$d = $word.ActiveDocument

# Find the paragraph that currently contains only "01/11/18" - that is both
# the insertion point and the paragraph whose trailing mark/bookmark need to
# be folded into the new final paragraph afterwards.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "01/11/18") {
        $target = $cand
    }
}
if ($target -eq $null) {
    throw "could not find the '01/11/18' paragraph to update"
}

$beforeCount = $d.Paragraphs.Count
$targetIndex = $target.Index

$newBodyXml = @'
<w:p>
  <w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>06/11/2018</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>Bríd:</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Created </w:t></w:r>
  <w:proofErr w:type="spellStart"/><w:r><w:t>InitialiseGame</w:t></w:r><w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> and updated main so that we can now run things cleanly from different classes using ‘</w:t></w:r>
  <w:proofErr w:type="spellStart"/><w:r><w:t>InitialiseGame</w:t></w:r><w:proofErr w:type="spellEnd"/>
  <w:r><w:t>’ methods</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Updated Player, Cards to allow the </w:t></w:r>
  <w:proofErr w:type="spellStart"/><w:r><w:t>InitialiseGame</w:t></w:r><w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> class to access them (I had to change a few to ‘static’ to make it easier to access them) and pass the correct things to them (small changes)</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Copied Lilianna’s Utility, Space, </w:t></w:r>
  <w:proofErr w:type="spellStart"/><w:r><w:t>BoardReader</w:t></w:r><w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t>Implemented first part of ‘moving’ a player in main</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">I couldn’t see the changes we made in the lab last </w:t></w:r>
  <w:proofErr w:type="gramStart"/><w:r><w:t>week..</w:t></w:r><w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> or the ‘spaces’ text file? I’m not sure where they are so I made temporary ones, just so I could keep going. We can delete these when we sort it out though!</w:t></w:r>
</w:p>
'@

$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  $newBodyXml +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the target Range's contents (the "01/11/18" run plus
# the bookmark pair) with the parsed fragment above.
$target.Range.InsertXML($xmlPayload)

# Word leaves the original paragraph mark behind as a trailing empty
# paragraph once its contents have been replaced/relocated this way (the
# mark that used to end the "01/11/18" paragraph). Merge it away so the new
# "...sort it out though!" paragraph is the document's true last paragraph,
# matching the target layout exactly.
$afterCount = $d.Paragraphs.Count
$added = $afterCount - $beforeCount
$lastNewIndex = $targetIndex + $added
$secondLast = $d.Paragraphs.Item($lastNewIndex - 1)
$lastPara = $d.Paragraphs.Item($lastNewIndex)
$mergeRange = $d.Range($secondLast.Range.End - 1, $lastPara.Range.End)
$mergeRange.Delete()

Write-Output "Done: paragraphs before=$beforeCount after=$($d.Paragraphs.Count)"
